$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text fields) ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Cells whose type/format flips between a number and the "N/A" placeholder text ---
# (copy donor cell with value+format in one shot to get the right OOXML style index and type)
$ws.Range("D14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("D14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))

$ws.Range("C15").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("C15").Copy($ws.Range("F22"))
$ws.Range("F22").Value = 1
$ws.Range("C15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 0
$ws.Range("C15").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("C15").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = 0
$ws.Range("C15").Copy($ws.Range("F29"))
$ws.Range("F29").Value = 1
$ws.Range("C15").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 1
$ws.Range("K14").Copy($ws.Range("H29"))
$ws.Range("H29").Value = 0
$ws.Range("C15").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$ws.Range("C15").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = 0
$ws.Range("C15").Copy($ws.Range("F30"))
$ws.Range("F30").Value = 1
$ws.Range("C15").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("K14").Copy($ws.Range("H30"))
$ws.Range("H30").Value = 0
$ws.Range("C15").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 1
$ws.Range("K14").Copy($ws.Range("E33"))
$ws.Range("E33").Value = -100
$ws.Range("C15").Copy($ws.Range("G33"))
$ws.Range("G33").Value = 1
$ws.Range("K14").Copy($ws.Range("H33"))
$ws.Range("H33").Value = -100
$ws.Range("C15").Copy($ws.Range("J33"))
$ws.Range("J33").Value = 1
$ws.Range("K14").Copy($ws.Range("K33"))
$ws.Range("K33").Value = -100

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("L14").Value = -75
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 9.090909090909
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -50
$ws.Range("J16").Value = 78
$ws.Range("K16").Value = -29.487179487179
$ws.Range("L16").Value = -26.666666666666
$ws.Range("M16").Value = -33.734939759036
$ws.Range("N16").Value = -82.866043613707
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -61.538461538461
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -35.714285714285
$ws.Range("I17").Value = 114
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = -5
$ws.Range("L17").Value = -1.724137931034
$ws.Range("M17").Value = 75.384615384615
$ws.Range("N17").Value = -20.279720279720
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 71
$ws.Range("J18").Value = 67
$ws.Range("K18").Value = 5.970149253731
$ws.Range("L18").Value = -18.390804597701
$ws.Range("M18").Value = 12.698412698412
$ws.Range("N18").Value = -87.863247863247
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -38.888888888888
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -48.529411764705
$ws.Range("I19").Value = 301
$ws.Range("J19").Value = 393
$ws.Range("K19").Value = -23.409669211195
$ws.Range("L19").Value = -14
$ws.Range("M19").Value = 69.101123595505
$ws.Range("N19").Value = 40
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -28.571428571428
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -29.166666666666
$ws.Range("I20").Value = 138
$ws.Range("J20").Value = 168
$ws.Range("K20").Value = -17.857142857142
$ws.Range("L20").Value = -19.767441860465
$ws.Range("M20").Value = 105.970149253731
$ws.Range("N20").Value = -85.381355932203
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -47.727272727272
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = -39.716312056737
$ws.Range("I21").Value = 692
$ws.Range("J21").Value = 831
$ws.Range("K21").Value = -16.726835138387
$ws.Range("L21").Value = -14.567901234567
$ws.Range("M21").Value = 49.137931034482
$ws.Range("N21").Value = -68.842863574966
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 27
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = -22.857142857142
$ws.Range("L23").Value = -47.058823529411
$ws.Range("M23").Value = 12.5
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 62.5
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = 26.041666666666
$ws.Range("I24").Value = 596
$ws.Range("J24").Value = 570
$ws.Range("K24").Value = 4.561403508771
$ws.Range("L24").Value = 10.166358595194
$ws.Range("M24").Value = 46.798029556650
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 11.111111111111
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 17.073170731707
$ws.Range("I25").Value = 177
$ws.Range("J25").Value = 240
$ws.Range("K25").Value = -26.25
$ws.Range("L25").Value = -20.982142857142
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = -10.344827586206
$ws.Range("I26").Value = 199
$ws.Range("J26").Value = 170
$ws.Range("K26").Value = 17.058823529411
$ws.Range("L26").Value = 24.375
$ws.Range("M26").Value = -2.926829268292
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 6.666666666666
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = -20.689655172413
$ws.Range("L28").Value = 64.285714285714
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 66.666666666666
$ws.Range("L29").Value = 66.666666666666
$ws.Range("M29").Value = 66.666666666666
$ws.Range("N29").Value = 25
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 33.333333333333
$ws.Range("L30").Value = 33.333333333333
$ws.Range("M30").Value = 33.333333333333
$ws.Range("N30").Value = 0
